# Työaikaraportti_Panu_Käppi.docx edit
#
# Commit message: "Muokattu classeja ja lisätty metodeita classeihin"
#   (Finnish: "Edited classes and added methods to classes")
#
# This fills in the first still-empty row of the time-tracking table with
# a new work-log entry: date, hours, and description.

$d = $word.ActiveDocument

# The table has a header row followed by data rows. The first empty data
# row (immediately after the "19.09.2022 / 2 / Classien tekeminen" entry)
# is row 5 (1-based: row 1 = header).
$table = $d.Tables.Item(1)
$row = $table.Rows.Item(5)

$row.Cells.Item(1).Range.Text = "20.09.2022"
$row.Cells.Item(2).Range.Text = "1"
$row.Cells.Item(3).Range.Text = "Metodien lisääminen classeihin ja classien muokkaus"
